$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update Row 2 (Real) values
$ws.Range("B2").Value = 21
$ws.Range("C2").Value = 21
$ws.Range("D2").Value = 17
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 13
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

# Update Row 3 (Estimado) base value; dependent formulas recalc automatically
$ws.Range("B3").Value = 21

# Update selection to I2
$ws.Range("I2").Select()
